# Split the single combined "group1 group2 meandiff p-adj lower upper reject"
# header string (A1) and the single combined data string (A2) out into
# separate columns A:G, matching the pandas Tukey HSD results dataframe.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1) - keep the existing bold/bordered style, applied
# across the full header range.
$ws.Range("A1").Value = "group1"
$ws.Range("B1").Value = "group2"
$ws.Range("C1").Value = "meandiff"
$ws.Range("D1").Value = "p-adj"
$ws.Range("E1").Value = "lower"
$ws.Range("F1").Value = "upper"
$ws.Range("G1").Value = "reject"

$ws.Range("A1").Copy()
$ws.Range("B1:G1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Data row (row 2)
$ws.Range("A2").Value = "Ambient"
$ws.Range("B2").Value = "Reduced"
$ws.Range("C2").Value = 0.1209
$ws.Range("D2").Value = 0.326
$ws.Range("E2").Value = -0.1232
$ws.Range("F2").Value = 0.3651
$ws.Range("G2").Value = $false
